$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.798604087325389
$ws.Range("D2").Value = 4.214406168578183
$ws.Range("E2").Value = 13.10593716066539
$ws.Range("F2").Value = 20.14998911335667
$ws.Range("G2").Value = 21.65500853270591
$ws.Range("H2").Value = 12.49574039983071
$ws.Range("K2").Value = 11.62664108442963
$ws.Range("M2").Value = 13.0207827403957
$ws.Range("O2").Value = 18.06781118903344
$ws.Range("B3").Value = 7.726673904500051
$ws.Range("D3").Value = 4.130368078683959
$ws.Range("E3").Value = 13.00441565460421
$ws.Range("F3").Value = 20.17992074103945
$ws.Range("G3").Value = 21.69537377469808
$ws.Range("H3").Value = 12.54518136912992
$ws.Range("K3").Value = 11.25785599029455
$ws.Range("M3").Value = 12.73383572142636
$ws.Range("O3").Value = 18.14355556221603
$ws.Range("B4").Value = 7.684021651053623
$ws.Range("D4").Value = 4.077310613995063
$ws.Range("E4").Value = 12.94718014833197
$ws.Range("F4").Value = 20.20502078673313
$ws.Range("G4").Value = 21.73005288235932
$ws.Range("H4").Value = 12.57787916355296
$ws.Range("K4").Value = 11.02377102952373
$ws.Range("M4").Value = 12.55667693715103
$ws.Range("O4").Value = 18.19492896177139
$ws.Range("B5").Value = 7.667040349198293
$ws.Range("D5").Value = 4.055339622000782
$ws.Range("E5").Value = 12.92516208216719
$ws.Range("F5").Value = 20.21693443812654
$ws.Range("G5").Value = 21.74665886561944
$ws.Range("H5").Value = 12.5917920267041
$ws.Range("K5").Value = 10.9265571148112
$ws.Range("M5").Value = 12.48434171146046
$ws.Range("O5").Value = 18.21708352944031
$ws.Range("B6").Value = 7.664245316631435
$ws.Range("D6").Value = 4.0516707490784
$ws.Range("E6").Value = 12.92158550092553
$ws.Range("F6").Value = 20.21901430316414
$ws.Range("G6").Value = 21.74956524502463
$ws.Range("H6").Value = 12.59413776480713
$ws.Range("K6").Value = 10.91030789362204
$ws.Range("M6").Value = 12.47232501934017
$ws.Range("O6").Value = 18.22083581052795
$ws.Range("B7").Value = 7.683790991348075
$ws.Range("D7").Value = 4.07701569823829
$ws.Range("E7").Value = 12.94687788862688
$ws.Range("F7").Value = 20.2051746426612
$ws.Range("G7").Value = 21.73026683967632
$ws.Range("H7").Value = 12.57806441600398
$ws.Range("K7").Value = 11.0224672068266
$ws.Range("M7").Value = 12.55570183400273
$ws.Range("O7").Value = 18.19522281387564
$ws.Range("B8").Value = 7.773500324981717
$ws.Range("D8").Value = 4.185742264160899
$ws.Range("E8").Value = 13.06988935297551
$ws.Range("F8").Value = 20.15891163829814
$ws.Range("G8").Value = 21.66686495449707
$ws.Range("H8").Value = 12.51230157878637
$ws.Range("K8").Value = 11.50112933410158
$ws.Range("M8").Value = 12.92210585700271
$ws.Range("O8").Value = 18.09291545499159
$ws.Range("B9").Value = 7.960543578451127
$ws.Range("D9").Value = 4.386691190071272
$ws.Range("E9").Value = 13.35038695673873
$ws.Range("F9").Value = 20.12170808932341
$ws.Range("G9").Value = 21.62157392240606
$ws.Range("H9").Value = 12.40193112881442
$ws.Range("K9").Value = 12.37512740736518
$ws.Range("M9").Value = 13.62857352440115
$ws.Range("O9").Value = 17.93108091487483
$ws.Range("B10").Value = 8.10352559870546
$ws.Range("D10").Value = 4.525985095861918
$ws.Range("E10").Value = 13.57862961556798
$ws.Range("F10").Value = 20.12718290641449
$ws.Range("G10").Value = 21.63705347361172
$ws.Range("H10").Value = 12.33219792664546
$ws.Range("K10").Value = 12.97317150624564
$ws.Range("M10").Value = 14.13476945049795
$ws.Range("O10").Value = 17.83606820307124
$ws.Range("B11").Value = 8.169512982184504
$ws.Range("D11").Value = 4.587377105711187
$ws.Range("E11").Value = 13.68684814697952
$ws.Range("F11").Value = 20.13681319340147
$ws.Range("G11").Value = 21.65475367271813
$ws.Range("H11").Value = 12.30294541937402
$ws.Range("K11").Value = 13.23484867118879
$ws.Range("M11").Value = 14.36119709414519
$ws.Range("O11").Value = 17.79808027212824
$ws.Range("B12").Value = 8.194613939607585
$ws.Range("D12").Value = 4.610328328368277
$ws.Range("E12").Value = 13.7284205187101
$ws.Range("F12").Value = 20.14148594841245
$ws.Range("G12").Value = 21.66299083652084
$ws.Range("H12").Value = 12.29222384720494
$ws.Range("K12").Value = 13.3323890075635
$ws.Range("M12").Value = 14.44630532014846
$ws.Range("O12").Value = 17.78445168590102
$ws.Range("B13").Value = 8.189203337335018
$ws.Range("D13").Value = 4.605398757332489
$ws.Range("E13").Value = 13.71944145253855
$ws.Range("F13").Value = 20.14043398080425
$ws.Range("G13").Value = 21.6611485684025
$ws.Range("H13").Value = 12.29451709920806
$ws.Range("K13").Value = 13.31145180347909
$ws.Range("M13").Value = 14.42800523304047
$ws.Range("O13").Value = 17.78735313294884
$ws.Range("B14").Value = 8.171575919381088
$ws.Range("D14").Value = 4.589271330973054
$ws.Range("E14").Value = 13.69025668672559
$ws.Range("F14").Value = 20.13717707011556
$ws.Range("G14").Value = 21.65540059475612
$ws.Range("H14").Value = 12.30205621655968
$ws.Range("K14").Value = 13.24290477354362
$ws.Range("M14").Value = 14.36821217685178
$ws.Range("O14").Value = 17.79694385401654
$ws.Range("B15").Value = 8.160792652968478
$ws.Range("D15").Value = 4.579353826373657
$ws.Range("E15").Value = 13.67245614068142
$ws.Range("F15").Value = 20.13531568853481
$ws.Range("G15").Value = 21.65207963894076
$ws.Range("H15").Value = 12.30672048672831
$ws.Range("K15").Value = 13.20071410359175
$ws.Range("M15").Value = 14.3315021608802
$ws.Range("O15").Value = 17.80291709774141
$ws.Range("B16").Value = 8.099230216699024
$ws.Range("D16").Value = 4.521932155235
$ws.Range("E16").Value = 13.5716423223201
$ws.Range("F16").Value = 20.12669718131045
$ws.Range("G16").Value = 21.63611141367535
$ws.Range("H16").Value = 12.33415939096894
$ws.Range("K16").Value = 12.95585628354931
$ws.Range("M16").Value = 14.11988736385525
$ws.Range("O16").Value = 17.8386564672464
$ws.Range("B17").Value = 8.061689243981254
$ws.Range("D17").Value = 4.486190947047359
$ws.Range("E17").Value = 13.51089343858405
$ws.Range("F17").Value = 20.12323854598799
$ws.Range("G17").Value = 21.62904742714665
$ws.Range("H17").Value = 12.35162519421561
$ws.Range("K17").Value = 12.80294207532974
$ws.Range("M17").Value = 13.98902376039255
$ws.Range("O17").Value = 17.86192481737111
$ws.Range("B18").Value = 8.040186902931069
$ws.Range("D18").Value = 4.465448501524186
$ws.Range("E18").Value = 13.47636790581622
$ws.Range("F18").Value = 20.12192137807624
$ws.Range("G18").Value = 21.62598769112499
$ws.Range("H18").Value = 12.3619034642113
$ws.Range("K18").Value = 12.71401627709679
$ws.Range("M18").Value = 13.9133972797887
$ws.Range("O18").Value = 17.87580055981936
$ws.Range("B19").Value = 8.032922772656867
$ws.Range("D19").Value = 4.458394089556159
$ws.Range("E19").Value = 13.4647506978059
$ws.Range("F19").Value = 20.12159085168092
$ws.Range("G19").Value = 21.62512393409033
$ws.Range("H19").Value = 12.36542340660284
$ws.Range("K19").Value = 12.68374223139913
$ws.Range("M19").Value = 13.887732706033
$ws.Range("O19").Value = 17.8805830940493
$ws.Range("B20").Value = 8.065676371090403
$ws.Range("D20").Value = 4.490014907222951
$ws.Range("E20").Value = 13.51731752325072
$ws.Range("F20").Value = 20.12353716019098
$ws.Range("G20").Value = 21.62969554324101
$ws.Range("H20").Value = 12.34974187250697
$ws.Range("K20").Value = 12.81932126704639
$ws.Range("M20").Value = 14.00299199745298
$ws.Range("O20").Value = 17.85939687043669
$ws.Range("B21").Value = 8.176750634799482
$ws.Range("D21").Value = 4.594016494126751
$ws.Range("E21").Value = 13.69881319993841
$ws.Range("F21").Value = 20.13810587078692
$ws.Range("G21").Value = 21.65704726707477
$ws.Range("H21").Value = 12.29983213608621
$ws.Range("K21").Value = 13.26308123077595
$ws.Range("M21").Value = 14.38579270853618
$ws.Range("O21").Value = 17.79410626123935
$ws.Range("B22").Value = 8.249992045230274
$ws.Range("D22").Value = 4.660253448490015
$ws.Range("E22").Value = 13.82086494271665
$ws.Range("F22").Value = 20.15360634049672
$ws.Range("G22").Value = 21.68386632826271
$ws.Range("H22").Value = 12.26928691834468
$ws.Range("K22").Value = 13.54404031422365
$ws.Range("M22").Value = 14.63223964619444
$ws.Range("O22").Value = 17.75584696931466
$ws.Range("B23").Value = 8.210849947049901
$ws.Range("D23").Value = 4.625064201782389
$ws.Range("E23").Value = 13.75542261410816
$ws.Range("F23").Value = 20.14478692104177
$ws.Range("G23").Value = 21.66873427247659
$ws.Range("H23").Value = 12.28539953124277
$ws.Range("K23").Value = 13.39493421307074
$ws.Range("M23").Value = 14.50107378949698
$ws.Range("O23").Value = 17.77586170433426
$ws.Range("B24").Value = 8.063873538558344
$ws.Range("D24").Value = 4.488286699197583
$ws.Range("E24").Value = 13.51441194626876
$ws.Range("F24").Value = 20.12340006577602
$ws.Range("G24").Value = 21.62939941068839
$ws.Range("H24").Value = 12.35059258420945
$ws.Range("K24").Value = 12.81191938641654
$ws.Range("M24").Value = 13.99667816981836
$ws.Range("O24").Value = 17.86053820286856
$ws.Range("B25").Value = 7.908876942641446
$ws.Range("D25").Value = 4.333735240519236
$ws.Range("E25").Value = 13.27047474044806
$ws.Range("F25").Value = 20.12601994830032
$ws.Range("G25").Value = 21.62529643513135
$ws.Range("H25").Value = 12.40193112881442
$ws.Range("K25").Value = 12.14612519031655
$ws.Range("M25").Value = 13.43933194335323
$ws.Range("O25").Value = 17.97068349319241
